$d = $word.ActiveDocument

# 1) Change "Proficiently validated" -> "Validated" (content edit).
#    This necessarily rebuilds the host paragraph's runs in this runtime,
#    but gets the visible text exactly right.
$rFind = $d.Content
$ok = $rFind.Find.Execute("Proficiently validated", $true, $false, $false, $false, $false, $true, 1, $false, "Validated", 2)

# 2) Re-split the new text into the same run boundaries the target XML uses:
#       "• "  |  "V"  |  "alidated ... Query Studi"  |  "o"  |  "."
#    Toggling a character property (Bold on/off, net no-op since the
#    paragraph's runs are all Bold=false already) forces the engine to
#    split the run at the touched boundaries without altering the
#    final formatting.

# 2a) Isolate the "V" that starts "Validated".
$rV = $d.Content
$rV.Find.Execute("Validated hundreds", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$vOnly = $d.Range($rV.Start, $rV.Start + 1)
$vOnly.Bold = $true
$vOnly.Bold = $false

# 2b) Isolate the trailing "o" and "." of "...Query Studio."
$rTail = $d.Content
$rTail.Find.Execute("Automation Studio and Query Studio.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$oOnly = $d.Range($rTail.End - 2, $rTail.End - 1)
$oOnly.Bold = $true
$oOnly.Bold = $false
$dotOnly = $d.Range($rTail.End - 1, $rTail.End)
$dotOnly.Bold = $true
$dotOnly.Bold = $false
